$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 191-193 with new values
$ws.Range("A191").Value = 58.7
$ws.Range("B191").Value = 9.5
$ws.Range("C191").Value = 4.4000000000000004

$ws.Range("A192").Value = 66.8
$ws.Range("B192").Value = 10.7
$ws.Range("C192").Value = 5.9

$ws.Range("A193").Value = 74.7
$ws.Range("B193").Value = 13.8
$ws.Range("C193").Value = 9.5

# Row 194 becomes a new data row (previously held the last values, now shifted down to 195)
$ws.Range("A194").Value = 81.5
$ws.Range("B194").Value = 16.8
$ws.Range("C194").Value = 11.9

# New row 195 holds what used to be in row 194 (with slight value correction) and uses style "2" (numFmt 164) like row 190
$ws.Range("A195").Value = 84.686999999999998
$ws.Range("B195").Value = 22.103719999999999
$ws.Range("C195").Value = 14.132149999999999

# Apply style: copy formatting from A190:C190 (style index "2") to A195:C195
$ws.Range("A190:C190").Copy()
$ws.Range("A195:C195").PasteSpecial(-4122) # xlPasteFormats

# Update the selection/view state to match target
$ws.Range("A195:C195").Select()
